$d = $word.ActiveDocument

# 1. Insert a new "November 24, 2024" paragraph right after the
#    "Module 6.2 Assignment: Movies: Setup" paragraph (same NoSpacing style).
$titleRng = $d.Content
$titleRng.Find.Execute("Module 6.2 Assignment: Movies: Setup", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$titleRng.Collapse(0)
$titleRng.InsertParagraphAfter()
$titleRng.Move(1, 1) | Out-Null
$titleRng.InsertAfter("November 24, 2024")

# 2. Mark the first screenshot's drawing run as NoProof (adds <w:rPr><w:noProof/></w:rPr>).
#    The drawing sits in the paragraph right after "Run the SQL script".
$sqlRng = $d.Content
$sqlRng.Find.Execute("Run the SQL script", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$sqlPara = $sqlRng.Paragraphs(1)
$firstImgPara = $sqlPara.Next()
$firstImgPara.Range.NoProofing = 1

# 3. Merge the split "S" / "how a list of database tables" runs into one run.
$d.Content.Find.Execute("Show a list of database tables", $true, $false, $false, $false, $false, $true, 1, $false, "Show a list of database tables", 2) | Out-Null

# 4. Mark the second screenshot's drawing run as NoProof.
#    The drawing sits two paragraphs after "Show a list of database tables"
#    (an empty paragraph sits in between).
$listRng = $d.Content
$listRng.Find.Execute("Show a list of database tables", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$listPara = $listRng.Paragraphs(1)
$secondImgPara = $listPara.Next().Next()
$secondImgPara.Range.NoProofing = 1

# 5. Merge the split "Run " / "mysql_test.py file " runs into one run.
$d.Content.Find.Execute("Run mysql_test.py file ", $true, $false, $false, $false, $false, $true, 1, $false, "Run mysql_test.py file ", 2) | Out-Null
